# "Fruta / hortaliza, semanal" — weekly refresh of the Acelga (Macroferia
# Regional de Talca) price sheet: a new week's record is inserted at the
# top of the data block (row 158) and every existing record for that
# market/product slides down one row to make room.
#
# Net effect on the data block (rows 158-279 before the edit):
#   - A brand-new row is inserted at row 158 (pushing old rows 158..279
#     down to 159..280, carrying their values/format with them).
#   - The new row 158 is populated with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 158, shifting rows 158:279 down
# to 159:280 (values + formatting move with them, e.g. D's date style).
$ws.Rows(158).Insert()

# Fill the newly inserted row 158 with the new week's data.
$ws.Range("A158").Value = 5
$ws.Range("B158").Value = "Macroferia Regional de Talca"
$ws.Range("C158").Value = "Maule"
$ws.Range("D158").Value = 44762
$ws.Range("E158").Value = 7
$ws.Range("F158").Value = 100112009
$ws.Range("G158").Value = "Acelga"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 300
$ws.Range("K158").Value = 5000
$ws.Range("L158").Value = 5000
$ws.Range("M158").Value = 5000
$ws.Range("N158").Value = "$/docena de atados (4 kilos)"
$ws.Range("O158").Value = "Región del Maule"
$ws.Range("P158").Value = 1250
$ws.Range("Q158").Value = 4
$ws.Range("R158").Value = "Hortaliza"
